# Update the division problems in the table to the new values.
# Each "old" text is unique within the document, and the replacements are
# ordered so that a value produced as a "new" replacement (e.g. 51÷4= and
# 76÷2=) is not itself re-matched by a later Find/Replace for the same text
# that appears earlier as an "old" value elsewhere in this list.

$d = $word.ActiveDocument

$d.Content.Find.Execute("94÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷4=", 2)
$d.Content.Find.Execute("70÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷3=", 2)
$d.Content.Find.Execute("11÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷4=", 2)
$d.Content.Find.Execute("51÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=", 2)
$d.Content.Find.Execute("51÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷2=", 2)
$d.Content.Find.Execute("93÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷3=", 2)
$d.Content.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷2=", 2)
$d.Content.Find.Execute("23÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷2=", 2)
$d.Content.Find.Execute("84÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷6=", 2)
$d.Content.Find.Execute("49÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=", 2)
$d.Content.Find.Execute("76÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷4=", 2)
$d.Content.Find.Execute("76÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=", 2)
$d.Content.Find.Execute("23÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷7=", 2)
$d.Content.Find.Execute("79÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=", 2)
$d.Content.Find.Execute("31÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=", 2)
$d.Content.Find.Execute("12÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=", 2)
$d.Content.Find.Execute("25÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷6=", 2)
$d.Content.Find.Execute("99÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷5=", 2)
$d.Content.Find.Execute("26÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=", 2)
$d.Content.Find.Execute("73÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=", 2)
$d.Content.Find.Execute("11÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷6=", 2)
$d.Content.Find.Execute("44÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=", 2)
$d.Content.Find.Execute("41÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=", 2)
$d.Content.Find.Execute("64÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=", 2)
$d.Content.Find.Execute("43÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=", 2)
